$wb = $excel.ActiveWorkbook

# --- "Calculations Rough Input" sheet: move the selection (no activation) ---
$wsCalc = $wb.Worksheets.Item("Calculations Rough Input")
$wsCalc.Range("B23:H23").Select()

# --- "EPS" sheet keeps its own selection (L11); it just stops being active ---
# (nothing to do here explicitly - activating TTC below clears tabSelected on EPS)

# --- "TTC" sheet: fill in the new hardware row (gimbal), move selection, and
#     make it the active sheet/tab ---
$wsTTC = $wb.Worksheets.Item("TTC")
$wsTTC.Activate()

$wsTTC.Range("A4").Value = "gimbal"
$wsTTC.Range("B4").Value = "y-"
$wsTTC.Range("C4").Value = "z+"
$wsTTC.Range("D4").Value = 0
$wsTTC.Range("E4").Value = 0
$wsTTC.Range("F4").Value = 45

$wsTTC.Range("C8").Select()

Write-Output "done"
